# Apply the edits described by the commit:
#  - Switch the example "input" row from an OSM-online setup to an OSM-dump
#    setup, renaming the test analysis and pointing it at the real pbf file.

$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("input")

# --- "input" sheet: row 2 example values -------------------------------
# network_source: Network based on OSM online -> Network based on OSM dump
$wsInput.Range("D2").Value = "Network based on OSM dump"

# OSM_area_of_interest (E2) is no longer used for the OSM-dump source;
# clear it and move the dump file name into name_of_pbf (F2) instead.
$wsInput.Range("E2").Value = ""
$wsInput.Range("F2").Value = "NL_with_margin_from_EU_dump.osm.pbf"

# analysis_name: test2 -> TestNL
$wsInput.Range("A2").Value = "TestNL"

# Restore the default view for the input sheet (no frozen/scrolled
# top-left cell, selection back on B7).
$wsInput.Activate()
$wsInput.Range("B7").Select()
